# Mac build instructions: "Test by typing wx-config-3 <en-dash>cflags ..."
# becomes "Test by typing wx-config-3 --cflags ..." (en dash -> two hyphens),
# and Word's hidden "last edit" bookmark (_GoBack) follows the edit, moving
# from its old spot (just before the 2nd "/local/wx-3.0.2") to right after
# "wx-config-3" (the point where the user actually typed).

$d = $word.ActiveDocument

# Locate the paragraph holding the "Test by typing wx-config-3 ..." bullet.
$targetRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text
    if ($paraText -like "*Test by typing wx-config-3*") {
        $targetRange = $d.Paragraphs($i).Range
        break
    }
}

if ($targetRange -ne $null) {
    # Scope the Find to this paragraph only, since the en dash character
    # also appears elsewhere in the document (e.g. "bunzip2 " / "make ").
    $foundDash = $targetRange.Find.Execute([char]0x2013, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($foundDash) {
        # targetRange now spans exactly the single en-dash character.
        $dashStart = $targetRange.Start
        $editPoint = $dashStart - 1   # right after "wx-config-3", before the space

        # Replace the en dash with two literal hyphens.
        $targetRange.Text = ""
        $d.Range($dashStart, $dashStart).InsertAfter("-")
        $d.Range($dashStart + 1, $dashStart + 1).InsertAfter("-")

        # Word drops the hidden _GoBack bookmark at the last edited spot;
        # re-seat it immediately after "wx-config-3".
        $editBookmarkRange = $d.Range($editPoint, $editPoint)
        $d.Bookmarks.Add("_GoBack", $editBookmarkRange)
    }
}
